$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "محمد"
$ws.Range("C19").Value = "مجاهد "
$ws.Range("D19").Value = 3333

$ws.Range("D19").Select() | Out-Null
